$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3426.3333
$ws.Range("I2").Value = 140
$ws.Range("K2").Value = 140
$ws.Range("M2").Value = -27
$ws.Range("H33").Value = 200.15
$ws.Range("I33").Value = 216.94118
$ws.Range("J33").Value = 105
$ws.Range("K33").Value = 216.94118
$ws.Range("L33").Value = 105
$ws.Range("M33").Value = 12.05882
$ws.Range("N33").Value = -563
$ws.Range("H40").Value = 5236.1816
$ws.Range("I40").Value = 3449.75
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 3449.75
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = -3274.75
$ws.Range("N40").Value = -10350
$ws.Range("H41").Value = 667.4286
$ws.Range("I41").Value = 599.25
$ws.Range("K41").Value = 599.25
$ws.Range("M41").Value = -159.25
$ws.Range("H53").Value = 275.88235
$ws.Range("I53").Value = 179.625
$ws.Range("J53").Value = 361.44446
$ws.Range("K53").Value = 179.625
$ws.Range("L53").Value = 361.44446
$ws.Range("M53").Value = 457.375
$ws.Range("N53").Value = -1635.44446
$ws.Range("H55").Value = 450.1
$ws.Range("I55").Value = 223.66667
$ws.Range("J55").Value = 789.75
$ws.Range("K55").Value = 223.66667
$ws.Range("L55").Value = 789.75
$ws.Range("M55").Value = -9.666670000000011
$ws.Range("N55").Value = -1217.75
$ws.Range("H98").Value = 3540.9
$ws.Range("I98").Value = 1712.1111
$ws.Range("K98").Value = 1712.1111
$ws.Range("M98").Value = -214.1111000000001
$ws.Range("H107").Value = 2026.6923
$ws.Range("I107").Value = 2026.6923
$ws.Range("K107").Value = 2026.6923
$ws.Range("M107").Value = -106.6922999999999
$ws.Range("H112").Value = 3078
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 3308.889
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 9926.667000000001
$ws.Range("N112").Value = -12142.667
$ws.Range("M112").Value = -1892
$ws.Range("H122").Value = 3540.9
$ws.Range("I122").Value = 1712.1111
$ws.Range("K122").Value = 5136.3333
$ws.Range("M122").Value = -2686.3333
$ws.Range("H127").Value = 932
$ws.Range("I127").Value = 932
$ws.Range("K127").Value = 2796
$ws.Range("M127").Value = 2164
$ws.Range("H132").Value = 1654.7858
$ws.Range("I132").Value = 1654.7858
$ws.Range("K132").Value = 4964.357400000001
$ws.Range("M132").Value = -2434.357400000001
$ws.Range("H137").Value = 2090.1875
$ws.Range("I137").Value = 1995.6923
$ws.Range("K137").Value = 5987.0769
$ws.Range("M137").Value = -3437.0769
$ws.Range("H138").Value = 7807.727
$ws.Range("J138").Value = 8088.5
$ws.Range("L138").Value = 24265.5
$ws.Range("N138").Value = -34545.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12149.895
$ws.Range("I32").Value = 9740.468999999999
$ws.Range("J32").Value = 25000.166
$ws.Range("K32").Value = 9740.468999999999
$ws.Range("L32").Value = 25000.166
$ws.Range("M32").Value = -9453.468999999999
$ws.Range("N32").Value = -25574.166
$ws.Range("H45").Value = 2990.5715
$ws.Range("I45").Value = 2990.5715
$ws.Range("K45").Value = 2990.5715
$ws.Range("M45").Value = -2613.5715
$ws.Range("H61").Value = 3839.4
$ws.Range("I61").Value = 3839.4
$ws.Range("K61").Value = 3839.4
$ws.Range("M61").Value = -3627.4
$ws.Range("H122").Value = 2692.9443
$ws.Range("I122").Value = 2652.875
$ws.Range("K122").Value = 7958.625
$ws.Range("M122").Value = -5508.625
$ws.Range("H132").Value = 1644.909
$ws.Range("I132").Value = 899.55554
$ws.Range("K132").Value = 2698.66662
$ws.Range("M132").Value = -168.66662
$ws.Range("H136").Value = 3839.4
$ws.Range("I136").Value = 3839.4
$ws.Range("K136").Value = 11518.2
$ws.Range("M136").Value = -8968.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2247.6
$ws.Range("I107").Value = 2157
$ws.Range("K107").Value = 2157
$ws.Range("M107").Value = -237

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 774
$ws.Range("J22").Value = 732
$ws.Range("L22").Value = 732
$ws.Range("N22").Value = -1432
$ws.Range("H41").Value = 25000
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H132").Value = 2210.3845
$ws.Range("I132").Value = 1615.8948
$ws.Range("J132").Value = 3824
$ws.Range("K132").Value = 4847.6844
$ws.Range("L132").Value = 11472
$ws.Range("M132").Value = -2317.6844
$ws.Range("N132").Value = -16532

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 20000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H34").Value = 1843.8462
$ws.Range("J34").Value = 1843.8462
$ws.Range("L34").Value = 5531.5386
$ws.Range("N34").Value = -5699.5386
$ws.Range("H39").Value = 16166.333
$ws.Range("J39").Value = 16166.333
$ws.Range("L39").Value = 48498.999
$ws.Range("N39").Value = -49086.999
$ws.Range("H55").Value = 14414
$ws.Range("J55").Value = 14414
$ws.Range("L55").Value = 43242
$ws.Range("N55").Value = -43596
$ws.Range("H98").Value = 2649.2856
$ws.Range("I98").Value = 2860.8333
$ws.Range("J98").Value = 1380
$ws.Range("K98").Value = 8582.499899999999
$ws.Range("L98").Value = 4140
$ws.Range("M98").Value = -7084.499899999999
$ws.Range("N98").Value = -7136

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1511
$ws.Range("J16").Value = 3430.3333
$ws.Range("L16").Value = 3430.3333
$ws.Range("N16").Value = -3770.3333
$ws.Range("H22").Value = 7459
$ws.Range("I22").Value = 5960
$ws.Range("J22").Value = 8583.25
$ws.Range("K22").Value = 5960
$ws.Range("L22").Value = 8583.25
$ws.Range("M22").Value = -5665
$ws.Range("N22").Value = -9173.25
$ws.Range("H27").Value = 7459
$ws.Range("I27").Value = 5960
$ws.Range("J27").Value = 8583.25
$ws.Range("K27").Value = 5960
$ws.Range("L27").Value = 8583.25
$ws.Range("M27").Value = -5853
$ws.Range("N27").Value = -8797.25
$ws.Range("H93").Value = 2321.2222
$ws.Range("I93").Value = 2065.3333
$ws.Range("J93").Value = 2833
$ws.Range("K93").Value = 2065.3333
$ws.Range("L93").Value = 2833
$ws.Range("M93").Value = -817.3332999999998
$ws.Range("N93").Value = -5329
$ws.Range("H136").Value = 5080.4
$ws.Range("I136").Value = 2996
$ws.Range("K136").Value = 8988
$ws.Range("M136").Value = -6438

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 90000
$ws.Range("J93").Value = 90000
$ws.Range("L93").Value = 90000
$ws.Range("N93").Value = -94992
$ws.Range("H136").Value = 983.82355
$ws.Range("I136").Value = 901.6667
$ws.Range("J136").Value = 1600
$ws.Range("K136").Value = 2705.0001
$ws.Range("L136").Value = 4800
$ws.Range("M136").Value = -155.0001000000002
$ws.Range("N136").Value = -9900
